$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster data (rows 7-16 reshuffled; Rui Hachimura replaced by Zach Collins)
$data = @(
    @(7,  "Jalen Johnson",      "PF",    "Atlanta Hawks"),
    @(8,  "Trey Murphy III",    "SF,PF", "New Orleans Pelicans"),
    @(9,  "Franz Wagner",       "SF,PF", "Orlando Magic"),
    @(10, "Daniel Gafford",     "PF,C",  "Dallas Mavericks"),
    @(11, "Jarrett Allen",      "C",     "Cleveland Cavaliers"),
    @(12, "Mark Williams",      "C",     "Charlotte Hornets"),
    @(13, "Zach Collins",       "PF,C",  "San Antonio Spurs"),
    @(14, "OG Anunoby",         "SF,PF", "New York Knicks"),
    @(15, "Karl-Anthony Towns", "PF,C",  "New York Knicks"),
    @(16, "Jalen Duren",        "C",     "Detroit Pistons")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
